# The add-on software stops scanning the diagnosedCasesList as soon as it
# finds a matching classification, so cases must be ordered from highest
# to lowest priority (top = highest priority). This adds a new
# "Lumbar disorder" / "FX Lumbar" case right after the other "Lumbar
# disorder" rows (and before "OA/RA/Gouty arthritis"), which pushes every
# row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 60 (shifts rows 60..90 down to 61..91) and fill
# in the new case.
$ws.Rows.Item(60).Insert()
$ws.Range("A60").Value = "Lumbar disorder"
$ws.Range("B60").Value = "FX Lumbar"

# Keep the named range "diagnosedCasesList" in sync with the new last row.
$wb.Names.Item("diagnosedCasesList").RefersTo = "=Sheet1!`$A`$1:`$B`$91"

# Reflect the editor's final scroll position / selection in the sheet view.
[void]$ws.Range("B66").Select()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
